$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the B2:D9 range to 0
$ws.Range("B2:D9").Value = 0

# Override the two cells that receive new non-zero values
$ws.Range("D4").Value = 0.6398124592709094
$ws.Range("D8").Value = 0.6324229047117447
